$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 23.65990166666667
$ws.Range("N2").Value = 70.979705
$ws.Range("O2").Value = 0.2997993941754699
$ws.Range("P2").Value = 0.29979939417547
$ws.Range("Q2").Value = 964.2212315298522
$ws.Range("R2").Value = 8677.99108376867
$ws.Range("S2").Value = 0.006375343937217643
$ws.Range("T2").Value = 0.006375343937217644

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07690566666666666
$ws.Range("N3").Value = 0.230717
$ws.Range("O3").Value = 0.0009744872400636476
$ws.Range("P3").Value = 0.0009744872400636479
$ws.Range("Q3").Value = 3.134166729417555
$ws.Range("R3").Value = 28.207500564758
$ws.Range("S3").Value = 0.00002072282812619527
$ws.Range("T3").Value = 0.00002072282812619528

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 53.21452433333334
$ws.Range("N4").Value = 159.643573
$ws.Range("O4").Value = 0.6742919890890982
$ws.Range("P4").Value = 0.6742919890890983
$ws.Range("Q4").Value = 2168.672334773522
$ws.Range("R4").Value = 19518.0510129617
$ws.Range("S4").Value = 0.01433906614913815
$ws.Range("T4").Value = 0.01433906614913816

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.967779333333333
$ws.Range("N5").Value = 5.903338
$ws.Range("O5").Value = 0.02493412949536815
$ws.Range("P5").Value = 0.02493412949536816
$ws.Range("Q5").Value = 80.19368122897912
$ws.Range("R5").Value = 721.7431310608121
$ws.Range("S5").Value = 0.0005302333973865702
$ws.Range("T5").Value = 0.0005302333973865703

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 23.65990166666667
$ws.Range("N6").Value = 70.979705
$ws.Range("O6").Value = 0.2997993941754699
$ws.Range("P6").Value = 0.29979939417547
$ws.Range("Q6").Value = 39968.41888228484
$ws.Range("R6").Value = 359715.7699405635
$ws.Range("S6").Value = 0.2642675857666601
$ws.Range("T6").Value = 0.2642675857666602

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07690566666666666
$ws.Range("N7").Value = 0.230717
$ws.Range("O7").Value = 0.0009744872400636476
$ws.Range("P7").Value = 0.0009744872400636479
$ws.Range("Q7").Value = 129.9159203220711
$ws.Range("R7").Value = 1169.24328289864
$ws.Range("S7").Value = 0.0008589923638781889
$ws.Range("T7").Value = 0.0008589923638781892

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 53.21452433333334
$ws.Range("N8").Value = 159.643573
$ws.Range("O8").Value = 0.6742919890890982
$ws.Range("P8").Value = 0.6742919890890983
$ws.Range("Q8").Value = 89894.72691565313
$ws.Range("R8").Value = 809052.5422408781
$ws.Range("S8").Value = 0.5943758377112663
$ws.Range("T8").Value = 0.5943758377112665

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.967779333333333
$ws.Range("N9").Value = 5.903338
$ws.Range("O9").Value = 0.02493412949536815
$ws.Range("P9").Value = 0.02493412949536816
$ws.Range("Q9").Value = 3324.148585679662
$ws.Range("R9").Value = 29917.33727111696
$ws.Range("S9").Value = 0.02197897104847905
$ws.Range("T9").Value = 0.02197897104847905

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 23.65990166666667
$ws.Range("N10").Value = 70.979705
$ws.Range("O10").Value = 0.2997993941754699
$ws.Range("P10").Value = 0.29979939417547
$ws.Range("Q10").Value = 2388.831490942104
$ws.Range("R10").Value = 21499.48341847893
$ws.Range("S10").Value = 0.01579473866038887
$ws.Range("T10").Value = 0.01579473866038887

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fn1"
$ws.Range("C11").Value = "Itga4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.07690566666666666
$ws.Range("N11").Value = 0.230717
$ws.Range("O11").Value = 0.0009744872400636476
$ws.Range("P11").Value = 0.0009744872400636479
$ws.Range("Q11").Value = 7.764811576713222
$ws.Range("R11").Value = 69.883304190419
$ws.Range("S11").Value = 0.00005134023478273034
$ws.Range("T11").Value = 0.00005134023478273036

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fn1"
$ws.Range("C12").Value = "Itga4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 53.21452433333334
$ws.Range("N12").Value = 159.643573
$ws.Range("O12").Value = 0.6742919890890982
$ws.Range("P12").Value = 0.6742919890890983
$ws.Range("Q12").Value = 5372.825859291957
$ws.Range("R12").Value = 48355.43273362761
$ws.Range("S12").Value = 0.03552464066095672
$ws.Range("T12").Value = 0.03552464066095672

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fn1"
$ws.Range("C13").Value = "Itga4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.967779333333333
$ws.Range("N13").Value = 5.903338
$ws.Range("O13").Value = 0.02493412949536815
$ws.Range("P13").Value = 0.02493412949536816
$ws.Range("Q13").Value = 198.6776320932184
$ws.Range("R13").Value = 1788.098688838966
$ws.Range("S13").Value = 0.001313638608866333
$ws.Range("T13").Value = 0.001313638608866333

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fn1"
$ws.Range("C14").Value = "Itga4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 23.65990166666667
$ws.Range("N14").Value = 70.979705
$ws.Range("O14").Value = 0.2997993941754699
$ws.Range("P14").Value = 0.29979939417547
$ws.Range("Q14").Value = 2020.857202986512
$ws.Range("R14").Value = 18187.71482687861
$ws.Range("S14").Value = 0.01336172581120331
$ws.Range("T14").Value = 0.01336172581120331

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fn1"
$ws.Range("C15").Value = "Itga4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.07690566666666666
$ws.Range("N15").Value = 0.230717
$ws.Range("O15").Value = 0.0009744872400636476
$ws.Range("P15").Value = 0.0009744872400636479
$ws.Range("Q15").Value = 6.568724275501556
$ws.Range("R15").Value = 59.118518479514
$ws.Range("S15").Value = 0.00004343181327653296
$ws.Range("T15").Value = 0.00004343181327653298

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fn1"
$ws.Range("C16").Value = "Itga4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 53.21452433333334
$ws.Range("N16").Value = 159.643573
$ws.Range("O16").Value = 0.6742919890890982
$ws.Range("P16").Value = 0.6742919890890983
$ws.Range("Q16").Value = 4545.198721346519
$ws.Range("R16").Value = 40906.78849211867
$ws.Range("S16").Value = 0.03005244456773692
$ws.Range("T16").Value = 0.03005244456773693

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fn1"
$ws.Range("C17").Value = "Itga4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.967779333333333
$ws.Range("N17").Value = 5.903338
$ws.Range("O17").Value = 0.02493412949536815
$ws.Range("P17").Value = 0.02493412949536816
$ws.Range("Q17").Value = 168.0734390057551
$ws.Range("R17").Value = 1512.660951051796
$ws.Range("S17").Value = 0.001111286440636197
$ws.Range("T17").Value = 0.001111286440636198

